$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''96.308.21'
$ws.Range('E2').Value = '  -1.40%  '
$ws.Range('D3').Value = '''3.318.99'
$ws.Range('E3').Value = '  -3.29%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '''248.38'
$ws.Range('E5').Value = '  -3.18%  '
$ws.Range('D6').Value = '''650.19'
$ws.Range('E6').Value = '  -1.24%  '
$ws.Range('D7').Value = '''1.38'
$ws.Range('E7').Value = '  -7.88%  '
$ws.Range('E8').Value = '  -3.32%  '
$ws.Range('D9').Value = '''0.999'
$ws.Range('E9').Value = '  +0.06%  '
$ws.Range('E10').Value = '  -7.81%  '
$ws.Range('D11').Value = '''3.316.90'
$ws.Range('E11').Value = '  -3.28%  '
$ws.Range('E12').Value = '  -3.99%  '
$ws.Range('D13').Value = '''39.98'
$ws.Range('E13').Value = '  -5.48%  '
$ws.Range('D14').Value = '''96.046.33'
$ws.Range('E14').Value = '  -1.39%  '
$ws.Range('D15').Value = '''6.05'
$ws.Range('E15').Value = '  -5.42%  '
$ws.Range('E16').Value = '  -4.92%  '
$ws.Range('D17').Value = '''3.936.05'
$ws.Range('E17').Value = '  -3.12%  '
$ws.Range('E18').Value = '  -3.51%  '
$ws.Range('D19').Value = '''3.319.11'
$ws.Range('E19').Value = '  -3.49%  '
$ws.Range('D20').Value = '''0.535'
$ws.Range('E20').Value = '  +2.36%  '
$ws.Range('D21').Value = '''16.99'
$ws.Range('E21').Value = '  -4.18%  '
$ws.Range('D22').Value = '''501.23'
$ws.Range('E22').Value = '  -1.88%  '
$ws.Range('D23').Value = '''10.43'
$ws.Range('E23').Value = '  -5.09%  '
$ws.Range('D24').Value = '''3.34'
$ws.Range('E24').Value = '  -4.07%  '
$ws.Range('E25').Value = '  -5.66%  '
$ws.Range('D26').Value = '''6.53'
$ws.Range('E26').Value = '  +4.68%  '
$ws.Range('D27').Value = '''95.71'
$ws.Range('E27').Value = '  -4.14%  '
$ws.Range('D28').Value = '''11.97'
$ws.Range('E28').Value = '  -7.27%  '
$ws.Range('E29').Value = '  -9.76%  '
$ws.Range('D30').Value = '''0.999'
$ws.Range('E30').Value = '  +0.02%  '
$ws.Range('D31').Value = '''10.98'
$ws.Range('E31').Value = '  -4.73%  '
$ws.Range('E32').Value = '  -7.04%  '
$ws.Range('D33').Value = '''2.45'
$ws.Range('E33').Value = '  +6.73%  '
$ws.Range('E34').Value = '  +0.05%  '
$ws.Range('D35').Value = '''0.542'
$ws.Range('E35').Value = '  -7.12%  '
$ws.Range('D36').Value = '''27.87'
$ws.Range('E36').Value = '  -7.61%  '
$ws.Range('E37').Value = '  +1.36%  '
$ws.Range('D38').Value = '''7.54'
$ws.Range('E38').Value = '  -4.47%  '
$ws.Range('E39').Value = '  -0.05%  '
$ws.Range('E40').Value = '  -3.74%  '
$ws.Range('D41').Value = '''503.00'
$ws.Range('E41').Value = '  -3.57%  '
$ws.Range('D43').Value = '''0.0430'
$ws.Range('E43').Value = '  +1.85%  '
$ws.Range('D44').Value = '''0.825'
$ws.Range('E44').Value = '  -5.47%  '
$ws.Range('E45').Value = '  -1.63%  '
$ws.Range('D46').Value = '''1.65'
$ws.Range('E46').Value = '  +3.45%  '
$ws.Range('E47').Value = '  -1.41%  '
$ws.Range('D48').Value = '''8.29'
$ws.Range('E48').Value = '  +0.10%  '
$ws.Range('D49').Value = '''52.99'
$ws.Range('E49').Value = '  +2.87%  '
$ws.Range('E50').Value = '  -6.83%  '
$ws.Range('D51').Value = '''161.39'
$ws.Range('E51').Value = '  -0.58%  '
